$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.7964389134426562
$ws.Range("J2").Value = 0.7964389134426563
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.774269
$ws.Range("N2").Value = 2.322807
$ws.Range("O2").Value = 0.02746805195987118
$ws.Range("P2").Value = 0.02746805195987118
$ws.Range("Q2").Value = 1.005700068817333
$ws.Range("R2").Value = 9.051300619356001
$ws.Range("S2").Value = 0.02187662545730622
$ws.Range("T2").Value = 0.02187662545730623

# Row 3
$ws.Range("I3").Value = 0.7964389134426562
$ws.Range("J3").Value = 0.7964389134426563
$ws.Range("M3").Value = 25.63013966666666
$ws.Range("N3").Value = 76.89041899999999
$ws.Range("O3").Value = 0.9092576457313354
$ws.Range("P3").Value = 0.9092576457313354
$ws.Range("Q3").Value = 33.29105676007244
$ws.Range("R3").Value = 299.619510840652
$ws.Range("S3").Value = 0.7241681714056923
$ws.Range("T3").Value = 0.7241681714056925

# Row 4
$ws.Range("I4").Value = 0.7964389134426562
$ws.Range("J4").Value = 0.7964389134426563
$ws.Range("M4").Value = 1.783575
$ws.Range("N4").Value = 5.350725000000001
$ws.Range("O4").Value = 0.06327430230879351
$ws.Range("P4").Value = 0.06327430230879351
$ws.Range("Q4").Value = 2.3166903237
$ws.Range("R4").Value = 20.85021291330001
$ws.Range("S4").Value = 0.05039411657965765
$ws.Range("T4").Value = 0.05039411657965766

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.3319853333333333
$ws.Range("H5").Value = 0.995956
$ws.Range("I5").Value = 0.2035610865573438
$ws.Range("J5").Value = 0.2035610865573438
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.774269
$ws.Range("N5").Value = 2.322807
$ws.Range("O5").Value = 0.02746805195987118
$ws.Range("P5").Value = 0.02746805195987118
$ws.Range("Q5").Value = 0.2570459520546666
$ws.Range("R5").Value = 2.313413568492
$ws.Range("S5").Value = 0.005591426502564953
$ws.Range("T5").Value = 0.005591426502564954

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.3319853333333333
$ws.Range("H6").Value = 0.995956
$ws.Range("I6").Value = 0.2035610865573438
$ws.Range("J6").Value = 0.2035610865573438
$ws.Range("M6").Value = 25.63013966666666
$ws.Range("N6").Value = 76.89041899999999
$ws.Range("O6").Value = 0.9092576457313354
$ws.Range("P6").Value = 0.9092576457313354
$ws.Range("Q6").Value = 8.508830460618221
$ws.Range("R6").Value = 76.57947414556399
$ws.Range("S6").Value = 0.185089474325643
$ws.Range("T6").Value = 0.185089474325643

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.3319853333333333
$ws.Range("H7").Value = 0.995956
$ws.Range("I7").Value = 0.2035610865573438
$ws.Range("J7").Value = 0.2035610865573438
$ws.Range("M7").Value = 1.783575
$ws.Range("N7").Value = 5.350725000000001
$ws.Range("O7").Value = 0.06327430230879351
$ws.Range("P7").Value = 0.06327430230879351
$ws.Range("Q7").Value = 0.5921207409
$ws.Range("R7").Value = 5.3290866681
$ws.Range("S7").Value = 0.01288018572913585
$ws.Range("T7").Value = 0.01288018572913585
